$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = 44904
$ws.Range("B8").Value = "Finalizing all the ACFs"
$ws.Range("C8").Value = 3

$ws.Range("B9").Value = "Write report"

$ws.Range("B10").Select()
